$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing volunteer names in row 2
$ws.Range("A2").Value = "FollowupVolunteer3"
$ws.Range("B2").Value = "FollowupVolunteer1"
$ws.Range("C2").Value = "FollowupVolunteer4"
$ws.Range("D2").Value = "FollowupVolunteer5"

# Add new row 3 with the new volunteer entry
$ws.Range("A3").Value = "FollowupVolunteer2"

# Update the selection to reflect the new data range
$ws.Range("A2:D3").Select()
